$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to make the edits, then re-protect at the end.
$ws.Unprotect()

# Update the confidential disclosure date in A42 (2021-04-26 -> 2021-04-27)
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-27 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for each holding row (2-39)
$ws.Range("D2").Value = 0.06307138765741144
$ws.Range("E2").Value = -0.002449524940617676
$ws.Range("D3").Value = 0.05695299619467755
$ws.Range("E3").Value = 0.001605811508315957
$ws.Range("D4").Value = 0.2919725671865165
$ws.Range("E4").Value = 0.02218741842860861
$ws.Range("D5").Value = 0.03711580271987302
$ws.Range("E5").Value = 0.002472865943091707
$ws.Range("D6").Value = 0.03260202463491996
$ws.Range("E6").Value = -0.001866805591734022
$ws.Range("D7").Value = 0.02950624659286403
$ws.Range("E7").Value = 0.004582890541976692
$ws.Range("D8").Value = 0.02858994681084452
$ws.Range("E8").Value = -0.005727516451377035
$ws.Range("D9").Value = 0.02402412603389939
$ws.Range("E9").Value = 0.003408019723007749
$ws.Range("D10").Value = 0.02514957646720923
$ws.Range("E10").Value = -0.008203711800790381
$ws.Range("D11").Value = 0.02309563214831688
$ws.Range("E11").Value = 0.001748944033790778
$ws.Range("D12").Value = 0.0223121609943943
$ws.Range("E12").Value = 0.01192590713017028
$ws.Range("D13").Value = 0.02206882330275771
$ws.Range("E13").Value = 0.00200792315623799
$ws.Range("D14").Value = 0.02143658080820123
$ws.Range("E14").Value = -0.0006094773731524228
$ws.Range("D15").Value = 0.02090352428101878
$ws.Range("E15").Value = 0.001875058595581125
$ws.Range("D16").Value = 0.02154981176985763
$ws.Range("E16").Value = -0.003435558025564722
$ws.Range("D17").Value = 0.02109307726586858
$ws.Range("E17").Value = 0.004129351949828219
$ws.Range("D18").Value = 0.01535411840060831
$ws.Range("E18").Value = -0.01344452008168828
$ws.Range("D19").Value = 0.01655959263855033
$ws.Range("E19").Value = -0.0003681885125185191
$ws.Range("D20").Value = 0.01549848787672022
$ws.Range("E20").Value = -0.01088865472427114
$ws.Range("D21").Value = 0.01576174986257136
$ws.Range("E21").Value = 0.01311063218390784
$ws.Range("D22").Value = 0.01607444151822251
$ws.Range("E22").Value = -0.04532646979138455
$ws.Range("D23").Value = 0.01518993350620652
$ws.Range("E23").Value = -0.001490868430860881
$ws.Range("D24").Value = 0.01447102577561303
$ws.Range("E24").Value = -0.004529278550631011
$ws.Range("D25").Value = 0.01405588187484778
$ws.Range("E25").Value = -0.005228505034856723
$ws.Range("D26").Value = 0.01468061193060205
$ws.Range("E26").Value = -0.00717898515255333
$ws.Range("D27").Value = 0.01289526451786964
$ws.Range("E27").Value = 0.003875379939209944
$ws.Range("D28").Value = 0.01326370064695163
$ws.Range("E28").Value = 0.01162332545311284
$ws.Range("D29").Value = 0.01431849060515089
$ws.Range("E29").Value = -0.005946225439503405
$ws.Range("D30").Value = 0.01297801022061855
$ws.Range("E30").Value = 0.009664429530201302
$ws.Range("D31").Value = 0.01248676204850923
$ws.Range("E31").Value = -0.003278459821428714
$ws.Range("D32").Value = 0.01334045817384371
$ws.Range("E32").Value = 0.0004488733279468526
$ws.Range("D33").Value = 0.01266172565945331
$ws.Range("E33").Value = -0.002192699600154779
$ws.Range("D34").Value = 0.006740726248145434
$ws.Range("E34").Value = -0.006218503682646426
$ws.Range("D35").Value = 0.005555938435890643
$ws.Range("E35").Value = -0.00930825004899083
$ws.Range("D36").Value = 0.00591675325024576
$ws.Range("E36").Value = -0.01067275136169599
$ws.Range("D37").Value = 0.005614731435212237
$ws.Range("E37").Value = 0.003141361256544517
$ws.Range("D38").Value = 0.005137310505535959
$ws.Range("E38").Value = -0.01165624668856624
$ws.Range("D39").Value = 0.9999999999999999
$ws.Range("E39").Value = 0.005495403421774459

# Restore sheet protection
$ws.Protect()

$wb.Save()
